$d = $word.ActiveDocument

# 1. Update "21 years" -> "15+ years" in the professional summary
$d.Content.Find.Execute(
    "Results-driven Marketing & Data Analytics Professional with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Results-driven Marketing & Data Analytics Professional with 15+ years of experience",
    2)

# 2. Rewrite the FLEEM bullet under RESEARCH DIRECTOR - Progressive Change Campaign Committee
$d.Content.Find.Execute(
    "Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2)

# 3. Add a new bullet point after "Developed innovative approaches..." (Lake Research Partners)
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Developed innovative approaches to visualizing demographic and market data*") {
        $targetIdx = $idx
    }
}

if ($targetIdx -gt 0) {
    $target = $d.Paragraphs.Item($targetIdx)
    $target.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.Text = "• Trained staff on building Python tooling for report generation and analysis"
}
